# Regenerate merged AHB files
# 1) Rename the "_old" / "_new" header-suffix labels in row 1 to
#    "_FV2410" / "_FV2504" respectively.
# 2) Turn the data range A1:U80 into an Excel Table ("Table1").
# 3) Freeze the header row (split/freeze pane below row 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Header renames -------------------------------------------------
$newHeadersFV2410 = @("Segmentname_FV2410","Segmentgruppe_FV2410","Segment_FV2410","Datenelement_FV2410","Segment ID_FV2410","Code_FV2410","Qualifier_FV2410","Beschreibung_FV2410","Bedingungsausdruck_FV2410","Bedingung_FV2410")
$newHeadersFV2504 = @("Segmentname_FV2504","Segmentgruppe_FV2504","Segment_FV2504","Datenelement_FV2504","Segment ID_FV2504","Code_FV2504","Qualifier_FV2504","Beschreibung_FV2504","Bedingungsausdruck_FV2504","Bedingung_FV2504")

# Columns A..J (1..10) held the "_old" labels
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $newHeadersFV2410[$i]
}

# Column K (11) is "diff" - unchanged

# Columns L..U (12..21) held the "_new" labels
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $newHeadersFV2504[$i]
}

# --- 2) Convert the range into an Excel Table ---------------------------
$range = $ws.Range("A1:U80")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $range, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# --- 3) Freeze the header row -------------------------------------------
$null = $ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
